# Apply vocabulary.xlsx update: new terms vocab:1263-1287 + renamed rows 279-282
# Date-like strings (column Y) are written with a leading apostrophe so the
# COM host keeps them as literal text instead of auto-converting to a date serial,
# matching the original inlineStr text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 50: G50 (append new vocab refs) and Y50 (date) ---
$ws.Range("G50").Value = "vocab:1000,vocab:1038,vocab:1034,vocab:1028,vocab:1064,vocab:1057,vocab:1252,vocab:1244,vocab:1263,vocab:1259"
$ws.Range("Y50").Value = "'2024-01-08"

# --- Rows 279-282 were renamed/repurposed; refresh label + date ---
$ws.Range("B279").Value = "urinary samples information"
$ws.Range("Y279").Value = "'2024-01-08"

$ws.Range("B280").Value = "creatinine"
$ws.Range("Y280").Value = "'2024-01-08"

$ws.Range("B281").Value = "specific gravity"
$ws.Range("Y281").Value = "'2024-01-08"

$ws.Range("B282").Value = "osmolarity"
$ws.Range("Y282").Value = "'2024-01-08"

# --- Append new rows 283-307 ---
$ws.Range("A283").Value = "vocab:1263"
$ws.Range("B283").Value = "blood samples information"
$ws.Range("Y283").Value = "'2024-01-08"

$ws.Range("A284").Value = "vocab:1264"
$ws.Range("B284").Value = "gravimetric"
$ws.Range("G284").Value = "vocab:1263"
$ws.Range("Y284").Value = "'2024-01-08"

$ws.Range("A285").Value = "vocab:1265"
$ws.Range("B285").Value = "enzymatic"
$ws.Range("G285").Value = "vocab:1263"
$ws.Range("Y285").Value = "'2024-01-08"

$ws.Range("A286").Value = "vocab:1266"
$ws.Range("B286").Value = "none"
$ws.Range("G286").Value = "vocab:1263,vocab:1259"
$ws.Range("Y286").Value = "'2024-01-08"

$ws.Range("A287").Value = "vocab:1267"
$ws.Range("B287").Value = "unit of measure"
$ws.Range("Y287").Value = "'2024-01-08"

$ws.Range("A288").Value = "vocab:1268"
$ws.Range("B288").Value = "mL"
$ws.Range("G288").Value = "vocab:1263"
$ws.Range("Y288").Value = "'2024-01-08"

$ws.Range("A289").Value = "vocab:1269"
$ws.Range("B289").Value = "µL"
$ws.Range("G289").Value = "vocab:1263"
$ws.Range("Y289").Value = "'2024-01-08"

$ws.Range("A290").Value = "vocab:1270"
$ws.Range("B290").Value = "g"
$ws.Range("G290").Value = "vocab:1263"
$ws.Range("Y290").Value = "'2024-01-08"

$ws.Range("A291").Value = "vocab:1271"
$ws.Range("B291").Value = "mg"
$ws.Range("G291").Value = "vocab:1263"
$ws.Range("Y291").Value = "'2024-01-08"

$ws.Range("A292").Value = "vocab:1272"
$ws.Range("B292").Value = "µg"
$ws.Range("G292").Value = "vocab:1263"
$ws.Range("Y292").Value = "'2024-01-08"

$ws.Range("A293").Value = "vocab:1273"
$ws.Range("B293").Value = "cm"
$ws.Range("G293").Value = "vocab:1263"
$ws.Range("Y293").Value = "'2024-01-08"

$ws.Range("A294").Value = "vocab:1274"
$ws.Range("B294").Value = "sampling container material"
$ws.Range("Y294").Value = "'2024-01-08"

$ws.Range("A295").Value = "vocab:1275"
$ws.Range("B295").Value = "PP"
$ws.Range("G295").Value = "vocab:1274"
$ws.Range("Y295").Value = "'2024-01-08"

$ws.Range("A296").Value = "vocab:1276"
$ws.Range("B296").Value = "glass"
$ws.Range("G296").Value = "vocab:1274"
$ws.Range("Y296").Value = "'2024-01-08"

$ws.Range("A297").Value = "vocab:1277"
$ws.Range("B297").Value = "PET"
$ws.Range("G297").Value = "vocab:1274"
$ws.Range("Y297").Value = "'2024-01-08"

$ws.Range("A298").Value = "vocab:1278"
$ws.Range("B298").Value = "quality assurance/qualtiy control method"
$ws.Range("Y298").Value = "'2024-01-08"

$ws.Range("A299").Value = "vocab:1279"
$ws.Range("B299").Value = "standard operating procedure"
$ws.Range("G299").Value = "vocab:1278"
$ws.Range("Y299").Value = "'2024-01-08"

$ws.Range("A300").Value = "vocab:1280"
$ws.Range("B300").Value = "trained fieldworkers"
$ws.Range("G300").Value = "vocab:1278"
$ws.Range("Y300").Value = "'2024-01-08"

$ws.Range("A301").Value = "vocab:1281"
$ws.Range("B301").Value = "control of background contamination in the sampling material"
$ws.Range("G301").Value = "vocab:1278"
$ws.Range("Y301").Value = "'2024-01-08"

$ws.Range("A302").Value = "vocab:1282"
$ws.Range("B302").Value = "controlf of the transprot conditions"
$ws.Range("G302").Value = "vocab:1278"
$ws.Range("Y302").Value = "'2024-01-08"

$ws.Range("A303").Value = "vocab:1283"
$ws.Range("B303").Value = "control of background contamination in the conservation material"
$ws.Range("G303").Value = "vocab:1278"
$ws.Range("Y303").Value = "'2024-01-08"

$ws.Range("A304").Value = "vocab:1284"
$ws.Range("B304").Value = "identifciation/traceability of the samples"
$ws.Range("G304").Value = "vocab:1278"
$ws.Range("Y304").Value = "'2024-01-08"

$ws.Range("A305").Value = "vocab:1285"
$ws.Range("B305").Value = "criteria for acceptation/acceptation of the samples"
$ws.Range("G305").Value = "vocab:1278"
$ws.Range("Y305").Value = "'2024-01-08"

$ws.Range("A306").Value = "vocab:1286"
$ws.Range("B306").Value = "collection of field blanks"
$ws.Range("G306").Value = "vocab:1278"
$ws.Range("Y306").Value = "'2024-01-08"

$ws.Range("A307").Value = "vocab:1287"
$ws.Range("B307").Value = "settings the conditions for sample storage"
$ws.Range("G307").Value = "vocab:1278"
$ws.Range("Y307").Value = "'2024-01-08"

